$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at row 77, shifting the existing rows 77-92 down
# to 78-93 (dimension grows from R92 to R93).
$ws.Rows(77).Insert()

# Populate the newly inserted row 77 with the new price record.
$ws.Cells.Item(77, 1).Value = 7
$ws.Cells.Item(77, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(77, 3).Value = "Ñuble"
$ws.Cells.Item(77, 4).Value = 44798
$ws.Cells.Item(77, 5).Value = 16
$ws.Cells.Item(77, 6).Value = 100112021
$ws.Cells.Item(77, 7).Value = "Ají"
$ws.Cells.Item(77, 8).Value = "Americana (o)"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 60
$ws.Cells.Item(77, 11).Value = 29000
$ws.Cells.Item(77, 12).Value = 30000
$ws.Cells.Item(77, 13).Value = 29500
$ws.Cells.Item(77, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(77, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(77, 16).Value = 1967
$ws.Cells.Item(77, 17).Value = 15
$ws.Cells.Item(77, 18).Value = "Hortaliza"
